# Insert a new weekly price record as row 136, shifting all the
# following rows (old 136..180) down by one (to 137..181).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 136..180 down to 137..181 by inserting a new row at 136.
# (Excel copies the formatting of the row above, which preserves the
# date-style (s="2") on column D for the new row.)
$ws.Rows.Item(136).Insert()

# Populate the newly inserted row 136 with the new weekly record.
$ws.Cells.Item(136, 1).Value  = 11
$ws.Cells.Item(136, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(136, 3).Value  = "Bíobío"
$ws.Cells.Item(136, 4).Value  = 44726
$ws.Cells.Item(136, 5).Value  = 8
$ws.Cells.Item(136, 6).Value  = 100112003
$ws.Cells.Item(136, 7).Value  = "Ajo"
$ws.Cells.Item(136, 8).Value  = "Chino"
$ws.Cells.Item(136, 9).Value  = "Primera"
$ws.Cells.Item(136, 10).Value = 310
$ws.Cells.Item(136, 11).Value = 16000
$ws.Cells.Item(136, 12).Value = 17000
$ws.Cells.Item(136, 13).Value = 16516
$ws.Cells.Item(136, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(136, 15).Value = "China"
$ws.Cells.Item(136, 16).Value = 1652
$ws.Cells.Item(136, 17).Value = 10
$ws.Cells.Item(136, 18).Value = "Hortaliza"
